$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a "1" (present) mark for săpt. 9 (column J) for the listed students,
# matching the exercise-tracking sheet for vector traversal exercises.
$rows = @(6, 9, 12, 13, 18, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Update the active selection to reflect where the user ended up (K19).
$ws.Range("K19").Select()
